$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-6: clear C and E cells entirely (they no longer exist in the sheet)
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("E6").ClearContents()

# Rows 7-19: update C and E values
$ws.Range("C7").Value = 0.9512119708358302
$ws.Range("E7").Value = 0.9990492459760025

$ws.Range("C8").Value = 1.149724574326472
$ws.Range("E8").Value = 1.044407816150583

$ws.Range("C9").Value = 1.525861534474027
$ws.Range("E9").Value = 1.137551461271413

$ws.Range("C10").Value = 1.634644186146694
$ws.Range("E10").Value = 1.2772981976928

$ws.Range("C11").Value = 1.518308876725216
$ws.Range("E11").Value = 1.265181861560016

$ws.Range("C12").Value = 1.593309007378396
$ws.Range("E12").Value = 1.33496666414632

$ws.Range("C13").Value = 1.565661119702044
$ws.Range("E13").Value = 1.412546132271975

$ws.Range("C14").Value = 0.287327989413555
$ws.Range("E14").Value = 0.9047322996724727

$ws.Range("C15").Value = -1.746350382706474
$ws.Range("E15").Value = 0.7112343933969312

$ws.Range("C16").Value = 5.778434165738466
$ws.Range("E16").Value = 1.531961367047852

$ws.Range("C17").Value = -0.2355225117835369
$ws.Range("E17").Value = 0.8847367780353999

$ws.Range("C18").Value = 0.1363842982220032
$ws.Range("E18").Value = 0.9899450936446508

$ws.Range("C19").Value = 0.678264046940269
$ws.Range("E19").Value = 1.007646955063968
